$d = $word.ActiveDocument

# --- Section 1: merge split runs back into single runs (also drops the
# proofErr spell/grammar-check markers Word had inserted around them). A
# same-text Find/Replace is enough to make Word re-flow the run list.

$d.Content.Find.Execute(
    "SSJAE-Formato Libre-DiagramasProcesos ", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "SSJAE-Formato Libre-DiagramasProcesos ", 2) | Out-Null

$d.Content.Find.Execute(
    "SSJAE-Especificacion de Requerimientos de Software ", $false, $false,
    $false, $false, $false, $true, 1, $false,
    "SSJAE-Especificacion de Requerimientos de Software ", 2) | Out-Null

$d.Content.Find.Execute(
    "SSJAE-DescripciónReglasIndicadoresNegocio", $false, $false, $false,
    $false, $false, $true, 1, $false,
    "SSJAE-DescripciónReglasIndicadoresNegocio", 2) | Out-Null

$d.Content.Find.Execute(
    "olivia.rodriguez", $false, $false, $false, $false, $false, $true, 1,
    $false, "olivia.rodriguez", 2) | Out-Null

$d.Content.Find.Execute(
    "Password:", $false, $false, $false, $false, $false, $true, 1, $false,
    "Password:", 2) | Out-Null

$d.Content.Find.Execute(
    "En cuanto al documento de arquitectura, se está trabajando en su elaboración, pero el desarrollador Misael Mendoza Antunez, nos puede apoyar en este tema.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "En cuanto al documento de arquitectura, se está trabajando en su elaboración, pero el desarrollador Misael Mendoza Antunez, nos puede apoyar en este tema.",
    2) | Out-Null

# --- Section 2: the "_GoBack" bookmark moves from its own trailing empty
# paragraph into the middle of the password, splitting "Senades2020" into
# "s" + bookmark + "enades2020" (note the leading capital becomes lower-case).

$pwd = $d.Paragraphs.Item(19)
$splitPos = $pwd.Range.Start + 1
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null

$firstChar = $d.Range($pwd.Range.Start, $pwd.Range.Start + 1)
$firstChar.Text = "s"
